$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.610.88"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "1.804.68"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.602"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.45"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.68%  "
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0671"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("E11").Value = "  +3.78%  "
$ws.Range("D12").Value = "2.067.59"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").Value = "1.798.09"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.634"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("D16").Value = "34.625.95"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("D20").Value = "0.0₃0768"
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.15%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0514"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.39%  "
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.641"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.305.73"
$ws.Range("E37").Value = "  -4.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0187"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "14.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.76%  "
$ws.Range("E41").Value = "  +4.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "82.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.941"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("E46").Value = "  +3.97%  "
$ws.Range("D47").Value = "1.965.81"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.49%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("E51").Value = "  -0.20%  "
